$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing existing rows 2-10 down to 3-11.
$ws.Rows("2:2").Insert()
# Excel's default row-insert copies formatting from the row above (the
# bold header row); strip that back to the plain/default style used by
# the rest of the data rows before filling in values.
$ws.Range("A2:R2").ClearFormats()

# Populate the newly inserted row 2 with the latest week's data.
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C2").Value = "Los Lagos"
$ws.Range("D2").Value = 44530
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 100112012
$ws.Range("G2").Value = "Espinaca"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 10000
$ws.Range("N2").Value = "`$/cuna 10 kilos"
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 1000
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = "Hortaliza"
